$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update patient record in row 13 with the new recommendation data
# nick (B13): typed through a couple of drafts before settling on the final name
$ws.Range("B13").Value = "ASDA"
$ws.Range("B13").Value = "ASDAD ASDASD"
$ws.Range("B13").Value = "LUIGGI STEEP"

# apellidos (C13): reuse existing "PASACHE LOPERA" value
$ws.Range("C13").Value = "PASACHE LOPERA"

# sexo (E13): update from M to F
$ws.Range("E13").Value = "F"

$wb.Save()
